$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 162
$ws1.Range("F5").Value = 1808
$ws1.Range("F8").Value = 159
$ws1.Range("F9").Value = 2239
$ws1.Range("F11").Value = 61
$ws1.Range("F13").Value = 1388
$ws1.Range("F14").Value = 488
$ws1.Range("F16").Value = 304
$ws1.Range("F23").Value = 56
$ws1.Range("F24").Value = 25
$ws1.Range("F25").Value = 1402
$ws1.Range("F30").Value = 349

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 162
$ws4.Range("F5").Value = 1808
$ws4.Range("F9").Value = 159
$ws4.Range("F10").Value = 2239
$ws4.Range("F12").Value = 61
$ws4.Range("F14").Value = 1388
$ws4.Range("F15").Value = 488
$ws4.Range("F17").Value = 304
$ws4.Range("F24").Value = 56
$ws4.Range("F25").Value = 25
$ws4.Range("F26").Value = 1402
$ws4.Range("F29").Value = 179
$ws4.Range("F31").Value = 349
